$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $NewValue) {
    $cell = $Worksheet.Range($Address)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = $origStyle
    "ok"
}

Set-TextValue $ws "D2" "64.036.56"
Set-TextValue $ws "E2" "  -1.19%  "

Set-TextValue $ws "D3" "3.097.44"
Set-TextValue $ws "E3" "  -1.78%  "

Set-TextValue $ws "E4" "  -0.53%  "

Set-TextValue $ws "D5" "598.64"
Set-TextValue $ws "E5" "  +1.04%  "

Set-TextValue $ws "D6" "156.73"
Set-TextValue $ws "E6" "  +2.31%  "

Set-TextValue $ws "D7" "1.00"
Set-TextValue $ws "E7" "  -0.26%  "

Set-TextValue $ws "E8" "  +0.52%  "

Set-TextValue $ws "D9" "3.094.30"
Set-TextValue $ws "E9" "  -1.76%  "

Set-TextValue $ws "D10" "0.158"
Set-TextValue $ws "E10" "  -2.39%  "

Set-TextValue $ws "E11" "  -0.84%  "

Set-TextValue $ws "E12" "  -3.01%  "

Set-TextValue $ws "D13" "0.0000239"
Set-TextValue $ws "E13" "  -3.69%  "

Set-TextValue $ws "D14" "36.92"
Set-TextValue $ws "E14" "  -4.47%  "

Set-TextValue $ws "E15" "  -0.82%  "

Set-TextValue $ws "D16" "3.607.59"
Set-TextValue $ws "E16" "  -1.85%  "

Set-TextValue $ws "D17" "7.21"
Set-TextValue $ws "E17" "  -1.38%  "

Set-TextValue $ws "D18" "63.947.01"
Set-TextValue $ws "E18" "  -0.74%  "

Set-TextValue $ws "D19" "3.101.37"
Set-TextValue $ws "E19" "  -1.77%  "

Set-TextValue $ws "D20" "481.28"
Set-TextValue $ws "E20" "  +1.21%  "

Set-TextValue $ws "D21" "14.48"
Set-TextValue $ws "E21" "  -3.50%  "

Set-TextValue $ws "D22" "0.714"
Set-TextValue $ws "E22" "  -4.90%  "

Set-TextValue $ws "E23" "  -1.56%  "

Set-TextValue $ws "E24" "  +3.38%  "

Set-TextValue $ws "D25" "81.59"
Set-TextValue $ws "E25" "  -0.94%  "

Set-TextValue $ws "D26" "12.91"
Set-TextValue $ws "E26" "  -4.53%  "

Set-TextValue $ws "D27" "10.80"
Set-TextValue $ws "E27" "  +8.68%  "

Set-TextValue $ws "E28" "  +0.00%  "

Set-TextValue $ws "E29" "  +2.82%  "

Set-TextValue $ws "E30" "  -1.48%  "

Set-TextValue $ws "B31" "FirstDigitalUSD"
Set-TextValue $ws "C31" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D31" "1.00"
Set-TextValue $ws "E31" "  -0.56%  "

Set-TextValue $ws "B32" "ImmutableX"
Set-TextValue $ws "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D32" "2.21"
Set-TextValue $ws "E32" "  -1.16%  "

Set-TextValue $ws "E33" "  -4.34%  "

Set-TextValue $ws "D34" "27.23"
Set-TextValue $ws "E34" "  -2.05%  "

Set-TextValue $ws "D35" "0.0₃0843"
Set-TextValue $ws "E35" "  -3.99%  "

Set-TextValue $ws "E36" "  +0.96%  "

Set-TextValue $ws "D37" "6.04"
Set-TextValue $ws "E37" "  -2.95%  "

Set-TextValue $ws "E38" "  -6.52%  "

Set-TextValue $ws "E39" "  -3.07%  "

Set-TextValue $ws "E40" "  -0.96%  "

Set-TextValue $ws "E41" "  -1.54%  "

Set-TextValue $ws "D42" "444.73"
Set-TextValue $ws "E42" "  -4.75%  "

Set-TextValue $ws "E43" "  -3.43%  "

Set-TextValue $ws "D44" "0.0365"
Set-TextValue $ws "E44" "  -4.45%  "

Set-TextValue $ws "E45" "  +1.09%  "

Set-TextValue $ws "B46" "Maker"
Set-TextValue $ws "C46" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D46" "2.840.24"
Set-TextValue $ws "E46" "  -2.08%  "

Set-TextValue $ws "B47" "Arweave"
Set-TextValue $ws "C47" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws "D47" "40.15"
Set-TextValue $ws "E47" "  +3.51%  "

Set-TextValue $ws "D48" "132.07"
Set-TextValue $ws "E48" "  +0.48%  "

Set-TextValue $ws "D49" "26.13"
Set-TextValue $ws "E49" "  +0.74%  "

Set-TextValue $ws "E50" "  -0.01%  "

Set-TextValue $ws "E51" "  -1.96%  "
